# Add two new rows (13 and 14) to Sheet1, continuing the existing table of
# "Getallen en variabelen" keyword/synonym pairs in columns A and B.
# Column A entries reuse the existing yellow-highlighted style used by the
# rest of the table (copied via Interior.Color, which the runtime maps back
# onto the existing cellXfs entry instead of creating a new one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: Kwadraat
$ws.Range("A13").Value = "Kwadraat"
$ws.Range("A13").Interior.Color = 65535
$ws.Range("B13").Value = "getal, vermenigvuldingen, Wortel"

# Row 14: Wortel
$ws.Range("A14").Value = "Wortel"
$ws.Range("A14").Interior.Color = 65535
$ws.Range("B14").Value = "getal, vermenigvuldingen,kwadraat, negatief"

# Match the final selection recorded in the saved workbook
$ws.Range("B14").Select() | Out-Null
